$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.574.82'
$ws.Range("E2").Value = '  -0.04%  '

$ws.Range("D3").Value = '2.370.56'
$ws.Range("E3").Value = '  -0.78%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").Value = "'508.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.49%  '

$ws.Range("D6").Value = "'133.78"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.27%  '

$ws.Range("D7").Value = "'0.996"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.35%  '

$ws.Range("E8").Value = '  -1.46%  '

$ws.Range("D9").Value = '2.392.46'
$ws.Range("E9").Value = '  -0.06%  '

$ws.Range("D10").Value = "'0.0974"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.03%  '

$ws.Range("E11").Value = '  -0.44%  '

$ws.Range("D12").Value = "'4.87"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.09%  '

$ws.Range("D13").Value = "'0.325"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.67%  '

$ws.Range("D14").Value = '2.794.84'
$ws.Range("E14").Value = '  -0.64%  '

$ws.Range("B15").Value = 'Avalanche'
$ws.Range("C15").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D15").Value = "'22.00"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.04%  '

$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '56.441.52'
$ws.Range("E16").Value = '  -0.20%  '

$ws.Range("E17").Value = '  +0.79%  '

$ws.Range("D18").Value = '2.372.24'
$ws.Range("E18").Value = '  +0.09%  '

$ws.Range("D19").Value = "'10.07"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.89%  '

$ws.Range("E20").Value = '  +1.01%  '

$ws.Range("D21").Value = "'312.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.71%  '

$ws.Range("D22").Value = "'6.30"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.75%  '

$ws.Range("E23").Value = '  -0.37%  '

$ws.Range("D24").Value = "'65.60"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.36%  '

$ws.Range("D25").Value = "'0.995"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.26%  '

$ws.Range("D26").Value = "'0.377"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.93%  '

$ws.Range("D27").Value = "'0.149"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.95%  '

$ws.Range("D28").Value = "'7.30"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.36%  '

$ws.Range("D29").Value = "'171.60"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.41%  '

$ws.Range("D30").Value = '0.0₃0722'
$ws.Range("E30").Value = '  -0.52%  '

$ws.Range("D31").Value = "'1.65"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.07%  '

$ws.Range("E32").Value = '  +0.52%  '

$ws.Range("E33").Value = '  -1.03%  '

$ws.Range("E34").Value = '  -0.07%  '

$ws.Range("E35").Value = '  -0.38%  '

$ws.Range("D36").Value = "'17.80"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.40%  '

$ws.Range("D37").Value = "'1.20"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.49%  '

$ws.Range("D38").Value = "'0.882"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +8.59%  '

$ws.Range("E39").Value = '  -1.62%  '

$ws.Range("D40").Value = "'36.55"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.32%  '

$ws.Range("E41").Value = '  +0.56%  '

$ws.Range("D42").Value = "'0.378"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.51%  '

$ws.Range("E43").Value = '  +4.34%  '

$ws.Range("E44").Value = '  +0.40%  '

$ws.Range("D45").Value = "'126.95"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.52%  '

$ws.Range("E46").Value = '  -0.69%  '

$ws.Range("D47").Value = "'0.0903"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.72%  '

$ws.Range("D48").Value = "'248.09"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.50%  '

$ws.Range("D49").Value = "'0.0488"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.70%  '

$ws.Range("E50").Value = '  +1.34%  '

$ws.Range("D51").Value = "'0.0211"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.76%  '
